$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter pairs (source -> destination) used to shift the comment
# text on row 15 one column to the left, for every commented column from
# U through BP (T itself is the column being deleted below).
$pairs = @(
    @("U15","T15"), @("V15","U15"), @("W15","V15"), @("X15","W15"),
    @("Y15","X15"), @("Z15","Y15"), @("AA15","Z15"), @("AB15","AA15"),
    @("AC15","AB15"), @("AD15","AC15"), @("AE15","AD15"), @("AF15","AE15"),
    @("AI15","AH15"), @("AJ15","AI15"), @("AK15","AJ15"),
    @("AM15","AL15"), @("AN15","AM15"), @("AO15","AN15"), @("AP15","AO15"),
    @("AQ15","AP15"), @("AR15","AQ15"), @("AS15","AR15"), @("AT15","AS15"),
    @("AU15","AT15"), @("AV15","AU15"), @("AW15","AV15"), @("AX15","AW15"),
    @("AY15","AX15"), @("AZ15","AY15"), @("BA15","AZ15"), @("BB15","BA15"),
    @("BC15","BB15"), @("BD15","BC15"), @("BE15","BD15"), @("BF15","BE15"),
    @("BG15","BF15"), @("BH15","BG15"), @("BI15","BH15"), @("BJ15","BI15"),
    @("BK15","BJ15"), @("BL15","BK15"), @("BM15","BL15"), @("BN15","BM15"),
    @("BO15","BN15"), @("BP15","BO15")
)

# Capture the existing comment text for every source cell before any
# writes happen (writes below would otherwise clobber values we still
# need to read for a later pair).
$texts = @{}
foreach ($pair in $pairs) {
    $src = $pair[0]
    $cmt = $ws.Range($src).Comment
    if ($cmt -eq $null) {
        $texts[$src] = $null
    } else {
        $texts[$src] = $cmt.Text()
    }
}

# Now delete column T entirely - this removes the "culture_collection"
# field (shared string + header cell) and shifts every later column
# left by one, matching what Excel does for Delete Column.
$ws.Range("T:T").EntireColumn.Delete()

# Re-home each comment's text onto its new column (comments are not
# moved automatically by EntireColumn.Delete in this host, so we must
# push the text across explicitly).
foreach ($pair in $pairs) {
    $src = $pair[0]
    $dst = $pair[1]
    $text = $texts[$src]
    $dstCmt = $ws.Range($dst).Comment
    if ($text -eq $null) {
        if ($dstCmt -ne $null) {
            $dstCmt.Delete()
        }
    } else {
        if ($dstCmt -eq $null) {
            $ws.Range($dst).AddComment($text)
        } else {
            $dstCmt.Text($text)
        }
    }
}

# The final shifted-from cell (old BP, now without a comment of its own
# since everything moved one step left) must no longer carry a comment.
$lastCmt = $ws.Range("BP15").Comment
if ($lastCmt -ne $null) {
    $lastCmt.Delete()
}
